$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$directions = "Combine the milk and breadcrumbs. Place the breadcrumbs in a small bowl, pour in the milk, and stir to combine. Set aside while preparing the rest of the meatball mixture. The breadcrumbs will absorb the milk and become soggy. /nWhisk the egg, salt, pepper, Parmesan, and parsley. Whisk the egg in a large bowl until broken up. Add the Parmesan, parsley, salt, and a generous quantity of black pepper, and whisk to combine. /nAdd the ground meat. Add the meat to the egg mixture. Use your hands to thoroughly mix the egg mixture into the ground meat. /nAdd the onions and soaked breadcrumbs. Add the onions, garlic, and soaked breadcrumbs. Mix them thoroughly into the meat with your fingers. Try not to overwork the meat; pinch the meat between your fingers rather than kneading it. /nForm the meat into meatballs. Have a rimmed baking sheet ready. Pinch off a piece of the meat mixture and gently roll between your hands to form 1 1/2-inch meatballs (about 2 tablespoons each). Continue shaping until all the meat is used, placing the meatballs on the baking sheet so that they are not touching. /nOption 1: Roast or broil the meatballs in the oven. Broil the meatballs on a rimmed baking sheet for 20 to 25 minutes, or roast at 400°F for 25 to 30 minutes. (Watch closely if making meatballs made with lean meat.) The meatballs are done when cooked through, the outsides are browned, and they register 165°F in the middle on an instant-read thermometer. Serve immediately. /nOption 2: Cook the meatballs directly in sauce. Bring a marinara or other pasta sauce to a simmer and add as many meatballs as will fit comfortably in the pan. Cover and simmer for 30 to 35 minutes. The meatballs are done when they're cooked through and register 165°F in the middle on an instant-read thermometer. Serve immediately."

$ingredients = "breadcrumbs milk egg parmesan-cheese parsley salt pepper beef onion garlic"

# Replace the old "meatballs" row (row 2) in place with the new "Classic Meatballs" recipe.
$ws.Range("C2").Value = $directions
$ws.Range("B2").Value = $ingredients
$ws.Range("A2").Value = "Classic Meatballs"

# Set the page to portrait orientation (page setup was touched in the edit).
$ws.PageSetup.Orientation = 1

# Move the selection/active cell to where the user last clicked.
$ws.Range("E16").Select()
